$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: "Pedro Paredes" -> "Maria Urbina"
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Maria"
$ws.Range("B2").Value = "Urbina"
$ws.Range("C2").Value = "maurbina"
$ws.Range("D2").Value = "murbina@yopmail.com"

# ---------------------------------------------------------------------------
# Row 3: "Mary James" -> "Emilia Urbina"
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Emilia"
$ws.Range("B3").Value = "Urbina"
$ws.Range("C3").Value = "eurbina2"
$ws.Range("D3").Value = "eurbina@yopmail.com"
$ws.Range("I3").Value = "Femenino"
$ws.Range("K3").Value = "Licenciado"
$ws.Range("L3").Value = "Genética aplicada"

# ---------------------------------------------------------------------------
# Row 4 (new): "Joan Magallanes"
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Joan"
$ws.Range("B4").Value = "Magallanes"
$ws.Range("C4").Value = "jmagallanes"
$ws.Range("D4").Value = "jmagallanes@yopmail.com"
$ws.Range("E4").Value = "Tecnología"
$ws.Range("F4").Value = "Computación"
$ws.Range("I4").Value = "Masculino"
$ws.Range("J4").Value = 156423651
$ws.Range("K4").Value = "Ingeniero"
$ws.Range("L4").Value = "Desarrollo de aplicaciones"
$ws.Range("N4").Value = 4125658574

# D4 uses the same "hyperlink" font style as D2/D3 (blue text)
$ws.Range("D4").Font.Color = 16711680

# ---------------------------------------------------------------------------
# Hyperlinks: drop every existing hyperlink (item-level Delete() is a no-op
# in this host, so clear them all in bulk) and re-add the two that remain:
# D2 -> murbina@yopmail.com and D4 -> jmagallanes@yopmail.com. D3 no longer
# carries a hyperlink.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:murbina@yopmail.com", "", "", "murbina@yopmail")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jmagallanes@yopmail.com", "", "", "jmagallanes@yopmail")

# Hyperlinks.Add() stamps its display text into the cell itself - restore
# the full e-mail addresses as the actual cell text (the hyperlink's own
# "display" property is left holding the shorter text set above).
$ws.Range("D2").Value = "murbina@yopmail.com"
$ws.Range("D4").Value = "jmagallanes@yopmail.com"

# ---------------------------------------------------------------------------
# Selection moved to D9 in the authored edit
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
